$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append after the last existing row (140)
$rows = @(
    @{ L = "Especial"; M = 310; N = 18000; O = 18000; P = 18000; S = 1200 },
    @{ L = "Primera";  M = 350; N = 15000; O = 15000; P = 15000; S = 1000 },
    @{ L = "Segunda";  M = 280; N = 12000; O = 12000; P = 12000; S = 800 }
)

$startRow = 141
$dateSerial = 44911

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = 9
    $ws.Cells.Item($r, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = $dateSerial
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(140, 4).NumberFormat
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100103
    $ws.Cells.Item($r, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($r, 9).Value = 100103003
    $ws.Cells.Item($r, 10).Value = "Damasco"
    $ws.Cells.Item($r, 11).Value = "Dina"
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = "$/caja 15 kilos granel"
    $ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = 15
}
